$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Remove the explicit centered paragraph alignment from the "Sets" / "Reps" /
# "Last Weight" sub-header cells (rows 5 and 10 of the single table), so the
# paragraphs fall back to the (unset) default alignment again.
foreach ($r in 5, 10) {
    foreach ($c in 2, 3, 4) {
        $cell = $t.Cell($r, $c)
        $cell.Range.ParagraphFormat.Alignment = 0
    }
}

# Strip the leftover personal "2x3,75" weight values from the Dumbbell
# Snatch / Sumo Deadlift rows, leaving the (now empty) paragraphs in place.
$d.Content.Find.Execute("2x3,75", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
